$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: Products task now has owners, progress, and description
$ws.Range("E4").Value = "Product class is almost done just have to conect with data working on Product Manager class"
$ws.Range("B4").Value = "Wajahat and Archibald"
$ws.Range("D4").Value = "In Progress"

# Row 7: Customers task - clear owners, progress, and description
$ws.Range("B7").Value = ""
$ws.Range("D7").Value = ""
$ws.Range("E7").Value = ""

# Row 8: Bills task - clear owners, progress, and description
$ws.Range("B8").Value = ""
$ws.Range("D8").Value = ""
$ws.Range("E8").Value = ""

# Update active cell selection
$ws.Range("B4").Select()
